$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.481462001800537
$ws.Range("B1").Value = 1.670943021774292
$ws.Range("C1").Value = 1.716025590896606
$ws.Range("D1").Value = 2.118695735931396
$ws.Range("E1").Value = 3.122179269790649
